# Apply the row permutation described by the commit diff.
# Rows 2,4,5,6,7,8,9,11,12,13 each take on the field values that
# (in the pre-edit workbook) belonged to a different row at the same
# site "rismon önö, Jmt" (rows 3 and 10 are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replaced with data from original row 4
$ws.Range("A2").Value = 111815512
$ws.Range("B2").Value = 56398
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = 458154.6107204149
$ws.Range("R2").Value = 7054646.336103803
$ws.Range("AC2").Value = "ringhack"

# Row 4: replaced with data from original row 8
$ws.Range("A4").Value = 111815513
$ws.Range("Q4").Value = 458173.7327805056
$ws.Range("R4").Value = 7054711.474791372
$ws.Range("AC4").Value = "ringhack gamla"

# Row 5: replaced with data from original row 12
$ws.Range("A5").Value = 111815510
$ws.Range("B5").Value = 56398
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("Q5").Value = 458203.7272220219
$ws.Range("R5").Value = 7054385.000644128
$ws.Range("AC5").Value = "ringhack"

# Row 6: replaced with data from original row 11
$ws.Range("A6").Value = 111815516
$ws.Range("B6").Value = 89423
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma"
$ws.Range("H6").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 458289.5512131723
$ws.Range("R6").Value = 7054475.069158822

# Row 7: replaced with data from original row 9
$ws.Range("A7").Value = 111815515
$ws.Range("B7").Value = 89423
$ws.Range("E7").Value = 5432
$ws.Range("F7").Value = "Granticka"
$ws.Range("G7").Value = "Porodaedalea chrysoloma"
$ws.Range("H7").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("Q7").Value = 458161.9437607233
$ws.Range("R7").Value = 7054459.400503729
$ws.Range("AC7").ClearContents()

# Row 8: replaced with data from original row 13
$ws.Range("A8").Value = 111815507
$ws.Range("Q8").Value = 458151.5539710881
$ws.Range("R8").Value = 7054482.225765129

# Row 9: replaced with data from original row 5
$ws.Range("A9").Value = 111815514
$ws.Range("Q9").Value = 458153.7808649908
$ws.Range("R9").Value = 7054482.19637617

# Row 11: replaced with data from original row 2
$ws.Range("A11").Value = 111815519
$ws.Range("B11").Value = 77515
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 458215.7474518137
$ws.Range("R11").Value = 7054621.063481365

# Row 12: replaced with data from original row 6
$ws.Range("A12").Value = 111815517
$ws.Range("B12").Value = 77515
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("Q12").Value = 458250.8216980004
$ws.Range("R12").Value = 7054375.482693202
$ws.Range("AC12").ClearContents()

# Row 13: replaced with data from original row 7
$ws.Range("A13").Value = 111815509
$ws.Range("Q13").Value = 458176.2590895323
$ws.Range("R13").Value = 7054362.673967168
